# Atualizado por script em 11-11-2023 20:45
#
# This script applies match-result/odds corrections to three pairs of
# existing rows (a 3-way rotation among rows 63-65 and a swap between
# rows 68-69) and appends four newly scraped matches (rows 123-126) to
# the "Sheet1" results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers: capture / write back the data columns (F..V) of a row.
# Columns A-E (index/country/tournament/season/date) are left untouched
# for rows that already exist; they already carry the correct values.
# ---------------------------------------------------------------------
function Get-RowValues($row) {
    $vals = @{}
    for ($c = 6; $c -le 22; $c++) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# ---- Rotate rows 63, 64, 65 ------------------------------------------
# new63 <- old65 (Villarreal x Girona)
# new64 <- old63 (Ath Bilbao x Getafe)
# new65 <- old64 (Real Madrid x Las Palmas)
$row63 = Get-RowValues 63
$row64 = Get-RowValues 64
$row65 = Get-RowValues 65

Set-RowValues 63 $row65
Set-RowValues 64 $row63
Set-RowValues 65 $row64

# ---- Swap rows 68, 69 -------------------------------------------------
# new68 <- old69 (Celta Vigo x Alaves)
# new69 <- old68 (Granada CF x Betis)
$row68 = Get-RowValues 68
$row69 = Get-RowValues 69

Set-RowValues 68 $row69
Set-RowValues 69 $row68

# ---------------------------------------------------------------------
# Append the four new matches as rows 123-126.
# ---------------------------------------------------------------------
function Add-MatchRow(
    [int]$Row,
    [int]$Indice,
    [double]$DataPartida,
    [string]$Home,
    [int]$HomeGols,
    [string]$Away,
    [int]$AwayGols,
    [double]$HomeOpen,
    [string]$HomeOpenDt,
    [double]$HomeClose,
    [string]$HomeCloseDt,
    [double]$DrawOpen,
    [string]$DrawOpenDt,
    [double]$DrawClose,
    [string]$DrawCloseDt,
    [double]$AwayOpen,
    [string]$AwayOpenDt,
    [double]$AwayClose,
    [string]$AwayCloseDt,
    [string]$Url
) {
    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = "spain"
    $ws.Cells.Item($Row, 3).Value = "laliga"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"
    $ws.Cells.Item($Row, 5).Value = $DataPartida
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpen
    $ws.Cells.Item($Row, 11).Value = $HomeOpenDt
    $ws.Cells.Item($Row, 12).Value = $HomeClose
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDt
    $ws.Cells.Item($Row, 14).Value = $DrawOpen
    $ws.Cells.Item($Row, 15).Value = $DrawOpenDt
    $ws.Cells.Item($Row, 16).Value = $DrawClose
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDt
    $ws.Cells.Item($Row, 18).Value = $AwayOpen
    $ws.Cells.Item($Row, 19).Value = $AwayOpenDt
    $ws.Cells.Item($Row, 20).Value = $AwayClose
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDt
    $ws.Cells.Item($Row, 22).Value = $Url
}

# Row, Indice, DataPartida, Home, HomeGols, Away, AwayGols,
# HomeOpen, HomeOpenDt, HomeClose, HomeCloseDt,
# DrawOpen, DrawOpenDt, DrawClose, DrawCloseDt,
# AwayOpen, AwayOpenDt, AwayClose, AwayCloseDt, Url
Add-MatchRow 123 122 45241.67708333334 `
    "Almeria" 1 "Real Sociedad" 3 `
    3.94 "29/10/2023 11:02" 4.38 "11/11/2023 16:14" `
    3.57 "29/10/2023 11:02" 3.8 "11/11/2023 16:14" `
    1.97 "29/10/2023 11:02" 1.85 "11/11/2023 16:13" `
    "https://www.betexplorer.com/football/spain/laliga/almeria-real-sociedad/GMRyAZcO/"

Add-MatchRow 124 123 45241.77083333334 `
    "Osasuna" 1 "Las Palmas" 1 `
    1.79 "29/10/2023 11:02" 1.97 "11/11/2023 18:12" `
    3.56 "29/10/2023 11:02" 3.36 "11/11/2023 18:28" `
    4.92 "29/10/2023 11:02" 4.47 "11/11/2023 18:28" `
    "https://www.betexplorer.com/football/spain/laliga/osasuna-las-palmas/UPiBOYCH/"

Add-MatchRow 125 124 45241.77083333334 `
    "Granada CF" 1 "Getafe" 1 `
    2.52 "29/10/2023 11:02" 2.49 "11/11/2023 18:27" `
    3.01 "29/10/2023 11:02" 3.19 "11/11/2023 18:23" `
    3.22 "29/10/2023 11:02" 3.19 "11/11/2023 18:27" `
    "https://www.betexplorer.com/football/spain/laliga/granada-cf-getafe/OSl3Qfr5/"

Add-MatchRow 126 125 45241.875 `
    "Real Madrid" 5 "Valencia" 1 `
    1.41 "29/10/2023 11:02" 1.32 "11/11/2023 20:58" `
    4.84 "29/10/2023 11:02" 5.79 "11/11/2023 20:59" `
    7.79 "29/10/2023 11:02" 10.11 "11/11/2023 20:59" `
    "https://www.betexplorer.com/football/spain/laliga/real-madrid-valencia/YwdjTdbn/"

# Apply the same visual styles used by the other data rows: bold/centered
# index column (A) and the datetime number format on the match-date
# column (E), by copying formats from the row directly above (122),
# which already carries the correct style for both columns.
$xlPasteFormats = -4122

$ws.Cells.Item(122, 1).Copy()
$ws.Cells.Item(123, 1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(124, 1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(125, 1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(126, 1).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(122, 5).Copy()
$ws.Cells.Item(123, 5).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(124, 5).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(125, 5).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(126, 5).PasteSpecial($xlPasteFormats)
